$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the localization "Status" text from "Ready for handoff" to
#    "In Translation" everywhere it occurs (Overview!E2:F2, zh-cn!C2,
#    de-de!C2).
# ------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.UsedRange.Replace("Ready for handoff", "In Translation") | Out-Null
}

# ------------------------------------------------------------------
# 2. Narrow the "Status" column(s) now that the text is shorter.
#    Target stored column width is ~13.41 characters; this engine's
#    ColumnWidth setter quantizes to the nearest 1/6 of a character,
#    so 12.5 is the input that lands closest to the target.
# ------------------------------------------------------------------
$newStatusWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = $newStatusWidth
$overview.Range("F1").ColumnWidth = $newStatusWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = $newStatusWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = $newStatusWidth
